$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Column B (Coin name) updates ---
$ws.Range("B9").Value = "USDC"
$ws.Range("B10").Value = "Cardano"

# --- Column C (Link) updates ---
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"

# --- Column D (Price) updates - forced text to avoid numeric auto-conversion ---
Set-TextValue $ws.Range("D2") "97.685.96"
Set-TextValue $ws.Range("D3") "3.369.25"
Set-TextValue $ws.Range("D5") "252.50"
Set-TextValue $ws.Range("D6") "657.68"
Set-TextValue $ws.Range("D7") "1.44"
Set-TextValue $ws.Range("D8") "0.421"
Set-TextValue $ws.Range("D9") "0.999"
Set-TextValue $ws.Range("D10") "1.04"
Set-TextValue $ws.Range("D11") "3.365.22"
Set-TextValue $ws.Range("D13") "41.35"
Set-TextValue $ws.Range("D14") "97.417.60"
Set-TextValue $ws.Range("D15") "6.07"
Set-TextValue $ws.Range("D16") "0.0000254"
Set-TextValue $ws.Range("D17") "3.989.44"
Set-TextValue $ws.Range("D18") "8.86"
Set-TextValue $ws.Range("D19") "3.366.60"
Set-TextValue $ws.Range("D20") "18.03"
Set-TextValue $ws.Range("D21") "0.529"
Set-TextValue $ws.Range("D22") "10.86"
Set-TextValue $ws.Range("D23") "508.71"
Set-TextValue $ws.Range("D25") "7.09"
Set-TextValue $ws.Range("D26") "0.0000199"
Set-TextValue $ws.Range("D27") "93.46"
Set-TextValue $ws.Range("D28") "12.27"
Set-TextValue $ws.Range("D29") "3.547.72"
Set-TextValue $ws.Range("D30") "11.32"
Set-TextValue $ws.Range("D31") "0.999"
Set-TextValue $ws.Range("D32") "0.141"
Set-TextValue $ws.Range("D34") "2.55"
Set-TextValue $ws.Range("D35") "0.997"
Set-TextValue $ws.Range("D36") "0.559"
Set-TextValue $ws.Range("D37") "28.61"
Set-TextValue $ws.Range("D38") "7.97"
Set-TextValue $ws.Range("D39") "1.49"
Set-TextValue $ws.Range("D40") "523.92"
Set-TextValue $ws.Range("D41") "0.151"
Set-TextValue $ws.Range("D44") "0.853"
Set-TextValue $ws.Range("D45") "0.0429"
Set-TextValue $ws.Range("D48") "3.66"
Set-TextValue $ws.Range("D49") "5.60"
Set-TextValue $ws.Range("D50") "55.95"
Set-TextValue $ws.Range("D51") "8.56"

# --- Column E (Volume 1h %) updates ---
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("E12").Value = "  -2.72%  "
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("E15").Value = "  -5.17%  "
$ws.Range("E16").Value = "  -4.62%  "
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("E21").Value = "  -10.11%  "
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("E25").Value = "  +9.94%  "
$ws.Range("E26").Value = "  -3.54%  "
$ws.Range("E27").Value = "  -7.09%  "
$ws.Range("E28").Value = "  -5.37%  "
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("E32").Value = "  -6.18%  "
$ws.Range("E33").Value = "  -5.48%  "
$ws.Range("E34").Value = "  +6.67%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("E37").Value = "  -4.77%  "
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("E47").Value = "  +8.30%  "
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("E49").Value = "  -5.94%  "
$ws.Range("E50").Value = "  +8.54%  "
$ws.Range("E51").Value = "  -6.54%  "

Write-Host "Applied cryptos update"
